function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws 2 4 '64.382.13'
Set-TextCell $ws 2 5 '  +1.45%  '

Set-TextCell $ws 3 4 '3.162.19'
Set-TextCell $ws 3 5 '  +2.16%  '

Set-TextCell $ws 4 5 '  +0.12%  '

Set-TextCell $ws 5 4 '593.17'
Set-TextCell $ws 5 5 '  +1.72%  '

Set-TextCell $ws 6 4 '147.17'
Set-TextCell $ws 6 5 '  +1.50%  '

Set-TextCell $ws 7 5 '  +0.03%  '

Set-TextCell $ws 8 4 '3.149.68'
Set-TextCell $ws 8 5 '  +2.00%  '

Set-TextCell $ws 9 4 '0.532'
Set-TextCell $ws 9 5 '  +0.91%  '

Set-TextCell $ws 10 4 '0.164'
Set-TextCell $ws 10 5 '  +3.65%  '

Set-TextCell $ws 11 4 '5.91'
Set-TextCell $ws 11 5 '  +4.65%  '

Set-TextCell $ws 12 4 '0.459'
Set-TextCell $ws 12 5 '  +0.71%  '

Set-TextCell $ws 13 5 '  +1.80%  '

Set-TextCell $ws 14 4 '37.46'
Set-TextCell $ws 14 5 '  +0.17%  '

Set-TextCell $ws 15 4 '3.676.44'
Set-TextCell $ws 15 5 '  +1.91%  '

Set-TextCell $ws 16 5 '  -0.06%  '

Set-TextCell $ws 17 4 '7.29'
Set-TextCell $ws 17 5 '  +2.71%  '

Set-TextCell $ws 18 4 '64.174.80'
Set-TextCell $ws 18 5 '  +1.37%  '

Set-TextCell $ws 19 4 '3.152.56'
Set-TextCell $ws 19 5 '  +1.92%  '

Set-TextCell $ws 20 4 '469.41'
Set-TextCell $ws 20 5 '  +2.12%  '

Set-TextCell $ws 21 4 '14.41'
Set-TextCell $ws 21 5 '  +1.25%  '

Set-TextCell $ws 22 4 '0.736'
Set-TextCell $ws 22 5 '  +1.55%  '

Set-TextCell $ws 23 5 '  +2.06%  '

Set-TextCell $ws 24 4 '2.39'
Set-TextCell $ws 24 5 '  +12.89%  '

Set-TextCell $ws 25 4 '13.21'
Set-TextCell $ws 25 5 '  +1.89%  '

Set-TextCell $ws 26 4 '81.45'
Set-TextCell $ws 26 5 '  +0.28%  '

Set-TextCell $ws 27 5 '  +0.04%  '

Set-TextCell $ws 28 4 '9.83'
Set-TextCell $ws 28 5 '  +11.43%  '

Set-TextCell $ws 29 5 '  +2.39%  '

Set-TextCell $ws 30 5 '  +1.68%  '

Set-TextCell $ws 31 4 '7.35'
Set-TextCell $ws 31 5 '  +8.46%  '

Set-TextCell $ws 32 5 '  +0.15%  '

Set-TextCell $ws 33 2 'EthereumClassic'
Set-TextCell $ws 33 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 33 4 '28.02'
Set-TextCell $ws 33 5 '  +4.94%  '

Set-TextCell $ws 34 2 'Hedera'
Set-TextCell $ws 34 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 34 4 '0.114'
Set-TextCell $ws 34 5 '  +7.08%  '

Set-TextCell $ws 35 4 '0.0₃0879'
Set-TextCell $ws 35 5 '  +3.72%  '

Set-TextCell $ws 36 4 '1.06'
Set-TextCell $ws 36 5 '  +3.31%  '

Set-TextCell $ws 37 4 '6.19'
Set-TextCell $ws 37 5 '  +2.99%  '

Set-TextCell $ws 38 5 '  -0.32%  '

Set-TextCell $ws 39 4 '3.24'
Set-TextCell $ws 39 5 '  -3.75%  '

Set-TextCell $ws 40 4 '469.84'
Set-TextCell $ws 40 5 '  +7.11%  '

Set-TextCell $ws 41 4 '9.47'
Set-TextCell $ws 41 5 '  +8.59%  '

Set-TextCell $ws 42 4 '51.38'
Set-TextCell $ws 42 5 '  +2.42%  '

Set-TextCell $ws 43 5 '  +9.44%  '

Set-TextCell $ws 44 5 '  +1.79%  '

Set-TextCell $ws 45 4 '2.912.10'
Set-TextCell $ws 45 5 '  +1.75%  '

Set-TextCell $ws 46 4 '39.83'
Set-TextCell $ws 46 5 '  +11.23%  '

Set-TextCell $ws 47 5 '  -0.44%  '

Set-TextCell $ws 48 4 '133.98'
Set-TextCell $ws 48 5 '  +8.55%  '

Set-TextCell $ws 50 2 'Stellar'
Set-TextCell $ws 50 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 50 4 '0.111'
Set-TextCell $ws 50 5 '  +1.13%  '

Set-TextCell $ws 51 2 'ThetaToken'
Set-TextCell $ws 51 3 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell $ws 51 4 '2.24'
Set-TextCell $ws 51 5 '  +4.69%  '
